$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.622546434402466
$ws.Range("B1").Value = 4.813064098358154
$ws.Range("C1").Value = 6.748772144317627
$ws.Range("D1").Value = 6.620687007904053
$ws.Range("E1").Value = 5.355889797210693
